$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.536.57"
$ws.Range("E2").Value = "  -0.14%  "
$ws.Range("D3").Value = "1.820.89"
$ws.Range("E3").Value = "  -0.44%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.35"
$ws.Range("E5").Value = "  -0.06%  "
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5146"
$ws.Range("E7").Value = "  -3.51%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3873"
$ws.Range("E8").Value = "  -2.76%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08438"
$ws.Range("E9").Value = "  +7.65%  "
$ws.Range("E10").Value = "  -0.45%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.112"
$ws.Range("E11").Value = "  -0.58%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.414"
$ws.Range("E12").Value = "  +1.09%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.07"
$ws.Range("E13").Value = "  +0.03%  "
$ws.Range("E14").Value = "  +0.00%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.513"
$ws.Range("E15").Value = "  -0.87%  "
$ws.Range("D16").Value = "1.813.17"
$ws.Range("E16").Value = "  -0.95%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001136"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "92.91"
$ws.Range("E18").Value = "  -0.41%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06694"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.80"
$ws.Range("E20").Value = "  -0.11%  "
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.098"
$ws.Range("E22").Value = "  -0.15%  "
$ws.Range("D23").Value = "28.563.61"
$ws.Range("E23").Value = "  -0.11%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.45"
$ws.Range("E24").Value = "  +2.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.276"
$ws.Range("E25").Value = "  +1.71%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "21.12"
$ws.Range("E26").Value = "  +1.30%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "159.07"
$ws.Range("E27").Value = "  +1.23%  "
$ws.Range("B28").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C28").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D28").Value = "2.024.00"
$ws.Range("E28").Value = "  -0.81%  "
$ws.Range("B29").Value = "LidoDAOToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.419"
$ws.Range("E29").Value = "  +0.13%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.07"
$ws.Range("E30").Value = "  +0.47%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.097"
$ws.Range("E31").Value = "  -4.46%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1081"
$ws.Range("E32").Value = "  -3.77%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.759"
$ws.Range("E33").Value = "  +0.22%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.07538"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.2230"
$ws.Range("E36").Value = "  -1.68%  "
$ws.Range("E37").Value = "  +0.71%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.207"
$ws.Range("E38").Value = "  -0.10%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.738"
$ws.Range("E39").Value = "  -2.87%  "
$ws.Range("E40").Value = "  +0.63%  "
$ws.Range("E41").Value = "  -1.23%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.192"
$ws.Range("E42").Value = "  -0.51%  "
$ws.Range("E43").Value = "  +0.53%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.59"
$ws.Range("E44").Value = "  +0.19%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.773"
$ws.Range("E45").Value = "  +1.53%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5937"
$ws.Range("E46").Value = "  +0.08%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "125.92"
$ws.Range("E47").Value = "  +0.26%  "
$ws.Range("E48").Value = "  -0.32%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.199"
$ws.Range("E49").Value = "  +0.42%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06988"
$ws.Range("E50").Value = "  +0.50%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "74.47"
$ws.Range("E51").Value = "  -0.24%  "
